$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 826, shifting existing rows 826:901 down to 827:902
$ws.Rows.Item(826).Insert()

# Populate the newly inserted row 826 with data
$ws.Cells.Item(826, 1).Value = 5
$ws.Cells.Item(826, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(826, 3).Value = "Maule"
$ws.Cells.Item(826, 4).Value = 45106
$ws.Cells.Item(826, 5).Value = 7
$ws.Cells.Item(826, 6).Value = "Fruta"
$ws.Cells.Item(826, 7).Value = 100102
$ws.Cells.Item(826, 8).Value = "Cítricos"
$ws.Cells.Item(826, 9).Value = 100102005
$ws.Cells.Item(826, 10).Value = "Naranja"
$ws.Cells.Item(826, 11).Value = "Fukumoto"
$ws.Cells.Item(826, 12).Value = "Primera"
$ws.Cells.Item(826, 13).Value = 700
$ws.Cells.Item(826, 14).Value = 7000
$ws.Cells.Item(826, 15).Value = 8000
$ws.Cells.Item(826, 16).Value = 7429
$ws.Cells.Item(826, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(826, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(826, 19).Value = 495
$ws.Cells.Item(826, 20).Value = 15
